$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (MG %) values added for the 46024 measurement block (rows 187-214)
$values = @{
    187 = 0.072
    188 = 0.093
    189 = 0.045
    190 = 0.041
    191 = 0.055
    192 = 0.088
    193 = 0.072
    194 = 0.037
    195 = 0.062
    196 = 0.086
    198 = 0.072
    199 = 0.058
    201 = 0.041
    202 = 0.072
    203 = 0.055
    204 = 0.075
    206 = 0.083
    208 = 0.093
    210 = 0.078
    211 = 0.072
    212 = 0.062
    213 = 0.037
    214 = 0.072
}

for ($row = 187; $row -le 214; $row++) {
    $cell = $ws.Range("D$row")
    if ($values.ContainsKey($row)) {
        $cell.Value = $values[$row]
    }
    $cell.NumberFormat = "0.0%"
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}

# Scroll / select to match the latest author view
$excel.ActiveWindow.ScrollRow = 184
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H197").Select() | Out-Null
